$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new species row for OFRA (2023 spp comparisons)
$ws.Range("A28").Value = "OFRA"
$ws.Range("B28").Value = "hermaphroditic"
$ws.Range("C28").Value = "broadcast"
$ws.Range("E28").Value = 4
$ws.Range("F28").Value = "A,B,C,D"

# Update selection to reflect the new last-used cell
$ws.Range("G28").Select()
